$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert a new column before the old "Portugal" column (EA),
# shifting Portugal..North Korea one column to the right.
$ws.Columns("EA:EA").Insert()

# Step 2: The old "North Korea" column (previously FU) got pushed to FV.
# Remove that displaced column so the sheet keeps its original A:FU extent.
$ws.Columns("FV:FV").Delete()

# Step 3: Label the newly inserted column as "North Korea".
$ws.Range("EA1").Value = "North Korea"

# Step 4: Rename the row-2 metric label.
$ws.Range("A2").Value = "net_gain_all_taxes_pc"

# Step 5: Write the new data row (one combined-tax value per country column).
$ws.Cells.Item(2, 2).Value = $null
$ws.Cells.Item(2, 3).Value = 143.409393886515
$ws.Cells.Item(2, 4).Value = 168.639558819865
$ws.Cells.Item(2, 5).Value = -1100.97513405339
$ws.Cells.Item(2, 6).Value = 6.25570550147051
$ws.Cells.Item(2, 7).Value = 152.758299538036
$ws.Cells.Item(2, 8).Value = -1289.9146083007
$ws.Cells.Item(2, 9).Value = -547.053671471371
$ws.Cells.Item(2, 10).Value = 146.951557918648
$ws.Cells.Item(2, 11).Value = 207.088311245121
$ws.Cells.Item(2, 12).Value = -275.176706448841
$ws.Cells.Item(2, 13).Value = 204.229050977001
$ws.Cells.Item(2, 14).Value = 206.46864253029
$ws.Cells.Item(2, 15).Value = 260.859226667799
$ws.Cells.Item(2, 16).Value = -7.98184572698369
$ws.Cells.Item(2, 17).Value = -431.792746475264
$ws.Cells.Item(2, 18).Value = -273.303735512854
$ws.Cells.Item(2, 19).Value = 167.996399987714
$ws.Cells.Item(2, 20).Value = 152.624378550863
$ws.Cells.Item(2, 21).Value = 131.30627605769
$ws.Cells.Item(2, 22).Value = 198.230011946533
$ws.Cells.Item(2, 23).Value = 59.4445327500726
$ws.Cells.Item(2, 24).Value = -106.492298242389
$ws.Cells.Item(2, 25).Value = -371.587532192645
$ws.Cells.Item(2, 26).Value = $null
$ws.Cells.Item(2, 27).Value = 134.324446810054
$ws.Cells.Item(2, 28).Value = 195.859887451739
$ws.Cells.Item(2, 29).Value = -558.549211353626
$ws.Cells.Item(2, 30).Value = -1800.62183091095
$ws.Cells.Item(2, 31).Value = -84.0049854705434
$ws.Cells.Item(2, 32).Value = -10.6778367752665
$ws.Cells.Item(2, 33).Value = 183.133990704133
$ws.Cells.Item(2, 34).Value = 195.660364897698
$ws.Cells.Item(2, 35).Value = 206.930361809342
$ws.Cells.Item(2, 36).Value = 126.086429908275
$ws.Cells.Item(2, 37).Value = 140.308494736565
$ws.Cells.Item(2, 38).Value = 215.31818437111
$ws.Cells.Item(2, 39).Value = 212.088613947147
$ws.Cells.Item(2, 40).Value = -39.9837424779363
$ws.Cells.Item(2, 41).Value = $null
$ws.Cells.Item(2, 42).Value = -244.882731674327
$ws.Cells.Item(2, 43).Value = -232.920617912088
$ws.Cells.Item(2, 44).Value = -511.438215541767
$ws.Cells.Item(2, 45).Value = 213.654391234779
$ws.Cells.Item(2, 46).Value = -1235.38854581829
$ws.Cells.Item(2, 47).Value = 93.5675998633401
$ws.Cells.Item(2, 48).Value = 170.077713417708
$ws.Cells.Item(2, 49).Value = 155.51960535635
$ws.Cells.Item(2, 50).Value = 196.092803195308
$ws.Cells.Item(2, 51).Value = 224.471045387976
$ws.Cells.Item(2, 52).Value = -230.246152472394
$ws.Cells.Item(2, 53).Value = -249.544431020225
$ws.Cells.Item(2, 54).Value = 218.766745173625
$ws.Cells.Item(2, 55).Value = -397.311912018252
$ws.Cells.Item(2, 56).Value = 39.5341913664369
$ws.Cells.Item(2, 57).Value = -610.361066380005
$ws.Cells.Item(2, 58).Value = -4.0258132605518
$ws.Cells.Item(2, 59).Value = -1385.00415909099
$ws.Cells.Item(2, 60).Value = 134.076255786894
$ws.Cells.Item(2, 61).Value = 200.678897196908
$ws.Cells.Item(2, 62).Value = 150.097292493941
$ws.Cells.Item(2, 63).Value = 207.817210836551
$ws.Cells.Item(2, 64).Value = 217.456104803615
$ws.Cells.Item(2, 65).Value = -45.8573113096867
$ws.Cells.Item(2, 66).Value = -229.797584749281
$ws.Cells.Item(2, 67).Value = 165.482294745936
$ws.Cells.Item(2, 68).Value = -94.5863298432474
$ws.Cells.Item(2, 69).Value = $null
$ws.Cells.Item(2, 70).Value = 217.511718475224
$ws.Cells.Item(2, 71).Value = -127.169910133841
$ws.Cells.Item(2, 72).Value = 238.433435388084
$ws.Cells.Item(2, 73).Value = -55.6106095555704
$ws.Cells.Item(2, 74).Value = 195.45238156266
$ws.Cells.Item(2, 75).Value = 224.407637557826
$ws.Cells.Item(2, 76).Value = -818.213639797646
$ws.Cells.Item(2, 77).Value = 188.900694090682
$ws.Cells.Item(2, 78).Value = 70.7718047161104
$ws.Cells.Item(2, 79).Value = -1861.18967704674
$ws.Cells.Item(2, 80).Value = -367.537086187999
$ws.Cells.Item(2, 81).Value = -233.736263016668
$ws.Cells.Item(2, 82).Value = 160.841316743479
$ws.Cells.Item(2, 83).Value = 169.162837243482
$ws.Cells.Item(2, 84).Value = -264.635924166443
$ws.Cells.Item(2, 85).Value = 2.43480496890075
$ws.Cells.Item(2, 86).Value = 211.455307149995
$ws.Cells.Item(2, 87).Value = 219.005838003631
$ws.Cells.Item(2, 88).Value = 239.468724533638
$ws.Cells.Item(2, 89).Value = -224.201161289987
$ws.Cells.Item(2, 90).Value = -574.278596844366
$ws.Cells.Item(2, 91).Value = 232.759302143842
$ws.Cells.Item(2, 92).Value = 216.797485697664
$ws.Cells.Item(2, 93).Value = 211.985505455395
$ws.Cells.Item(2, 94).Value = 105.150138933505
$ws.Cells.Item(2, 95).Value = 214.091285005871
$ws.Cells.Item(2, 96).Value = 241.321010947748
$ws.Cells.Item(2, 97).Value = -209.652183461914
$ws.Cells.Item(2, 98).Value = -3386.3416322839
$ws.Cells.Item(2, 99).Value = -217.361056427962
$ws.Cells.Item(2, 100).Value = 204.356992684075
$ws.Cells.Item(2, 101).Value = 186.462831948461
$ws.Cells.Item(2, 102).Value = 234.133175290046
$ws.Cells.Item(2, 103).Value = 75.6966013677964
$ws.Cells.Item(2, 104).Value = 24.5515982731845
$ws.Cells.Item(2, 105).Value = 169.117449122951
$ws.Cells.Item(2, 106).Value = 192.698302979678
$ws.Cells.Item(2, 107).Value = -595.98411888471
$ws.Cells.Item(2, 108).Value = 268.803662829531
$ws.Cells.Item(2, 109).Value = 85.3270220511585
$ws.Cells.Item(2, 110).Value = 165.435265264414
$ws.Cells.Item(2, 111).Value = 203.525425910764
$ws.Cells.Item(2, 112).Value = 132.901993880669
$ws.Cells.Item(2, 113).Value = -81.5724400104323
$ws.Cells.Item(2, 114).Value = 213.734161298425
$ws.Cells.Item(2, 115).Value = -21.0091960911196
$ws.Cells.Item(2, 116).Value = 165.122118692209
$ws.Cells.Item(2, 117).Value = 193.164631427579
$ws.Cells.Item(2, 118).Value = 195.497146770179
$ws.Cells.Item(2, 119).Value = 228.213887238075
$ws.Cells.Item(2, 120).Value = -601.887478358581
$ws.Cells.Item(2, 121).Value = -655.277447730575
$ws.Cells.Item(2, 122).Value = 255.169719676921
$ws.Cells.Item(2, 123).Value = -534.939188033917
$ws.Cells.Item(2, 124).Value = -355.556693015911
$ws.Cells.Item(2, 125).Value = 227.982504358712
$ws.Cells.Item(2, 126).Value = -199.308480741729
$ws.Cells.Item(2, 127).Value = 116.971412145985
$ws.Cells.Item(2, 128).Value = 204.065018550283
$ws.Cells.Item(2, 129).Value = 197.759222385163
$ws.Cells.Item(2, 130).Value = -4.73464896585228
$ws.Cells.Item(2, 131).Value = $null
$ws.Cells.Item(2, 132).Value = -144.042952432229
$ws.Cells.Item(2, 133).Value = 163.626138363074
$ws.Cells.Item(2, 134).Value = $null
$ws.Cells.Item(2, 135).Value = 57.0895819099619
$ws.Cells.Item(2, 136).Value = -46.5139714935541
$ws.Cells.Item(2, 137).Value = 224.703848455585
$ws.Cells.Item(2, 138).Value = -238.491821535875
$ws.Cells.Item(2, 139).Value = 196.712120573614
$ws.Cells.Item(2, 140).Value = 199.809777691515
$ws.Cells.Item(2, 141).Value = -3915.52272967454
$ws.Cells.Item(2, 142).Value = 185.811234915227
$ws.Cells.Item(2, 143).Value = -1564.80103092644
$ws.Cells.Item(2, 144).Value = 153.205812495172
$ws.Cells.Item(2, 145).Value = 196.620951951001
$ws.Cells.Item(2, 146).Value = 118.756331000429
$ws.Cells.Item(2, 147).Value = $null
$ws.Cells.Item(2, 148).Value = 170.726207926267
$ws.Cells.Item(2, 149).Value = -163.325595930992
$ws.Cells.Item(2, 150).Value = -241.605540378305
$ws.Cells.Item(2, 151).Value = -810.107862367249
$ws.Cells.Item(2, 152).Value = 182.726001759872
$ws.Cells.Item(2, 153).Value = $null
$ws.Cells.Item(2, 154).Value = 192.84362075354
$ws.Cells.Item(2, 155).Value = 214.812653603885
$ws.Cells.Item(2, 156).Value = 124.988431111307
$ws.Cells.Item(2, 157).Value = 228.052001494068
$ws.Cells.Item(2, 158).Value = 93.4142640597317
$ws.Cells.Item(2, 159).Value = 223.011586520781
$ws.Cells.Item(2, 160).Value = -289.490171104686
$ws.Cells.Item(2, 161).Value = 206.733308827961
$ws.Cells.Item(2, 162).Value = 69.3449893348494
$ws.Cells.Item(2, 163).Value = $null
$ws.Cells.Item(2, 164).Value = 203.283037943749
$ws.Cells.Item(2, 165).Value = 200.389667143952
$ws.Cells.Item(2, 166).Value = 213.458735771433
$ws.Cells.Item(2, 167).Value = -153.3522359672
$ws.Cells.Item(2, 168).Value = -1256.84772859951
$ws.Cells.Item(2, 169).Value = 226.239722625196
$ws.Cells.Item(2, 170).Value = 191.266724746965
$ws.Cells.Item(2, 171).Value = 245.062183400689
$ws.Cells.Item(2, 172).Value = 138.808538433483
$ws.Cells.Item(2, 173).Value = 173.883739737477
$ws.Cells.Item(2, 174).Value = 225.520312798273
$ws.Cells.Item(2, 175).Value = 143.300669973747
$ws.Cells.Item(2, 176).Value = 199.045278066756
$ws.Cells.Item(2, 177).Value = 204.938181987822
